# Insert a new row for "Dead Space 2" above the "Cuphead" row (row 56),
# pushing the existing rows 56-68 down to 57-69.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 56
$lastRow = 68

for ($r = $lastRow; $r -ge $firstRow; $r--) {
    $ws.Cells.Item($r + 1, 1).Value = $ws.Cells.Item($r, 1).Value2

    $bVal = $ws.Cells.Item($r, 2).Value2
    if ($bVal -eq $null) {
        $ws.Cells.Item($r + 1, 2).Clear()
    } else {
        $ws.Cells.Item($r + 1, 2).Value = $bVal
    }

    $ws.Cells.Item($r + 1, 3).Value = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r + 1, 4).Value = $ws.Cells.Item($r, 4).Value2
}

# Write the new "Dead Space 2" row into the now-vacated row 56.
$ws.Cells.Item($firstRow, 1).Value = "Dead Space 2"
$ws.Cells.Item($firstRow, 2).Clear()
$ws.Cells.Item($firstRow, 3).Value = 9
$ws.Cells.Item($firstRow, 4).Value = 7.2
